{"js": "// The commit adds a note that the particle effects show in the point of\n// contact for the falling ball. In the TODO list this item is now done,\n// so the \"Make the particle effect show in the hitting point\" TODO entry\n// is removed, and a new sentence is appended to the 1/8/2017 changelog\n// paragraph describing the finished feature.\n\n// 1) Remove the TODO paragraph that is now complete.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === \"Make the particle effect show in the hitting point\") {\n    paragraphs.items[i].delete();\n    break;\n  }\n}\nawait context.sync();\n\n// 2) Append the new sentence to the 1/8/2017 changelog paragraph, right\n//    after \"... Made the spell seal closer.\" as its own run.\nconst datedParagraphs = context.document.body.paragraphs;\ndatedParagraphs.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < datedParagraphs.items.length; i++) {\n  if (datedParagraphs.items[i].text.indexOf(\"1/8/2017\") === 0) {\n    const endRange = datedParagraphs.items[i].getRange(\"End\");\n    endRange.insertText(\" The particle effects show in the point of contact.\", \"End\");\n    break;\n  }\n}\nawait context.sync();\n", "ps1": "# The commit adds a note that the particle effects show in the point of\n# contact for the falling ball. In the TODO list this item is now done,\n# so the \"Make the particle effect show in the hitting point\" TODO entry\n# is removed, and a new sentence is appended to the 1/8/2017 changelog\n# paragraph describing the finished feature.\n\n$d = $word.ActiveDocument\n\n# 1) Remove the TODO paragraph that is now complete.\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.TrimEnd(\"`r\") -eq \"Make the particle effect show in the hitting point\") {\n        $p.Range.Delete()\n        break\n    }\n}\n\n# 2) Append the new sentence to the 1/8/2017 changelog paragraph, right\n#    after \"... Made the spell seal closer.\" as its own run.\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$found = $rng.Find.Execute(\"Made the spell seal closer.\")\nif ($found) {\n    $rng.Collapse(0)\n    $rng.InsertAfter(\" The particle effects show in the point of contact.\")\n}\n"}
